$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 39
$ws.Range("I2").Value = 124
$ws.Range("J2").Value = 430
$ws.Range("K2").Value = 6
$ws.Range("L2").Value = 105
$ws.Range("M2").Value = 9
$ws.Range("N2").Value = 61
$ws.Range("O2").Value = 1
$ws.Range("R2").Value = 4
$ws.Range("S2").Value = 46
$ws.Range("T2").Value = 79
$ws.Range("U2").Value = 8
$ws.Range("V2").Value = 609
$ws.Range("X2").Value = 690
$ws.Range("Y2").Value = 2
$ws.Range("Z2").Value = 8
$ws.Range("AA2").Value = 4
